$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet to reflect the new "through" date
$ws.Name = "Through 2021-12-19"

# Update the label for the December row
$ws.Range("A13").Value = "December (through 12-19)"

# Update December row (row 13) values
$ws.Range("B13").Value = 24
$ws.Range("C13").Value = 60
$ws.Range("D13").Value = 74
$ws.Range("F13").Value = 30
$ws.Range("G13").Value = 89
$ws.Range("H13").Value = 133

# Update Total row (row 14) values
$ws.Range("B14").Value = 315
$ws.Range("C14").Value = 623
$ws.Range("D14").Value = 895
$ws.Range("F14").Value = 564
$ws.Range("G14").Value = 1353
$ws.Range("H14").Value = 1776
